$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the notes for 44497 (D61) and 44498 (D62) rows to reflect progress
# through lesson 3.2 / 3.3, and bump the hours logged on 44498 (C62).
$ws.Range("D61").Value = "Finished 3 small problems, 15 mins on 5 problems from 3.2"
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = "Finished 5 problems from 3.2, 2 problems from 3.3"

# Move the active selection to C62 to match the saved cursor position.
$ws.Range("C62").Select()
